# Update cryptos list (prices & 1h volume %) as published by the
# "Updated cryptos list ... with GitHub Actions" workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Assign a value while forcing it to remain plain text, even when it
    # looks like a number (e.g. "8.40" or "141.21"). A leading apostrophe
    # forces Excel to treat the input as text; resetting the style back to
    # Normal afterwards avoids leaving a stray "quote prefix" cell style.
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range('D2').Value = '63.288.28'
$ws.Range('E2').Value = '  -1.13%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '3.233.42'
$ws.Range('E3').Value = '  +2.78%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.02%  '

# Row 5 - BNB
Set-TextValue 'D5' '594.61'
$ws.Range('E5').Value = '  -1.32%  '

# Row 6 - Solana
Set-TextValue 'D6' '141.21'
$ws.Range('E6').Value = '  -1.55%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  +0.08%  '

# Row 8 - LidoStakedEther
$ws.Range('D8').Value = '3.228.09'
$ws.Range('E8').Value = '  +2.71%  '

# Row 9 - XRP
$ws.Range('E9').Value = '  -1.66%  '

# Row 10 - Dogecoin
$ws.Range('E10').Value = '  -1.40%  '

# Row 11 - Toncoin
Set-TextValue 'D11' '5.36'
$ws.Range('E11').Value = '  -0.77%  '

# Row 12 - Cardano
$ws.Range('E12').Value = '  -0.79%  '

# Row 13 - ShibaInu
$ws.Range('E13').Value = '  -3.05%  '

# Row 14 - Avalanche
Set-TextValue 'D14' '34.37'
$ws.Range('E14').Value = '  -2.03%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '3.760.77'
$ws.Range('E15').Value = '  +2.64%  '

# Row 16 - TRON
$ws.Range('E16').Value = '  +0.29%  '

# Row 17 - WrappedEther
$ws.Range('D17').Value = '3.225.59'
$ws.Range('E17').Value = '  +2.88%  '

# Row 18 - WrappedBTC
$ws.Range('D18').Value = '63.306.46'

# Row 19 - Polkadot
Set-TextValue 'D19' '6.78'
$ws.Range('E19').Value = '  -1.26%  '

# Row 20 - BitcoinCash
Set-TextValue 'D20' '474.02'
$ws.Range('E20').Value = '  -3.02%  '

# Row 21 - Chainlink
Set-TextValue 'D21' '14.19'
$ws.Range('E21').Value = '  -3.59%  '

# Row 22 - Polygon
Set-TextValue 'D22' '0.729'
$ws.Range('E22').Value = '  +2.27%  '

# Row 23 - Uniswap
$ws.Range('E23').Value = '  +2.11%  '

# Row 24 - Litecoin
Set-TextValue 'D24' '83.82'
$ws.Range('E24').Value = '  -5.01%  '

# Row 25 - InternetComputer(DFINITY)
Set-TextValue 'D25' '13.16'
$ws.Range('E25').Value = '  -1.08%  '

# Row 26 - Dai
$ws.Range('E26').Value = '  -0.10%  '

# Row 27 - NEARProtocol
Set-TextValue 'D27' '7.59'
$ws.Range('E27').Value = '  +8.34%  '

# Row 28 - PancakeSwap
$ws.Range('E28').Value = '  -1.10%  '

# Row 29 - RenderToken
Set-TextValue 'D29' '8.10'
$ws.Range('E29').Value = '  -1.32%  '

# Row 30 - ImmutableX
Set-TextValue 'D30' '2.12'
$ws.Range('E30').Value = '  +2.31%  '

# Row 31 - EthereumClassic
Set-TextValue 'D31' '27.39'
$ws.Range('E31').Value = '  -1.21%  '

# Row 32 - FirstDigitalUSD
$ws.Range('E32').Value = '  -0.08%  '

# Row 33 - Hedera
$ws.Range('E33').Value = '  -4.31%  '

# Row 34 - Stacks
$ws.Range('E34').Value = '  -4.84%  '

# Row 35 - Mantle
$ws.Range('E35').Value = '  -1.72%  '

# Row 36 - Filecoin
Set-TextValue 'D36' '5.92'
$ws.Range('E36').Value = '  -2.53%  '

# Row 37 - OKB
Set-TextValue 'D37' '52.64'

# Row 38 - PEPE
$ws.Range('E38').Value = '  -5.76%  '

# Row 39 - VeChain
$ws.Range('E39').Value = '  -1.45%  '

# Row 40 - Bittensor
Set-TextValue 'D40' '422.12'
$ws.Range('E40').Value = '  -2.57%  '

# Row 41 - Cosmos
Set-TextValue 'D41' '8.40'
$ws.Range('E41').Value = '  +0.08%  '

# Row 42 - Maker
$ws.Range('D42').Value = '2.968.40'
$ws.Range('E42').Value = '  +1.05%  '

# Row 43 - dogwifhat
Set-TextValue 'D43' '2.75'
$ws.Range('E43').Value = '  -6.92%  '

# Row 44 - Kaspa
$ws.Range('E44').Value = '  -9.42%  '

# Row 45 - TheGraph
Set-TextValue 'D45' '0.268'
$ws.Range('E45').Value = '  +2.52%  '

# Row 46 - Fetch.AI
Set-TextValue 'D46' '2.17'
$ws.Range('E46').Value = '  -1.68%  '

# Row 47 / 48 - ThetaToken and USDe swapped ranking positions
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D47' '0.999'
$ws.Range('E47').Value = '  +0.05%  '

$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D48' '2.36'
$ws.Range('E48').Value = '  -2.09%  '

# Row 49 - InjectiveProtocol
Set-TextValue 'D49' '25.98'
$ws.Range('E49').Value = '  +0.22%  '

# Row 50 - Stellar
$ws.Range('E50').Value = '  -0.58%  '

# Row 51 - Monero
Set-TextValue 'D51' '120.90'
$ws.Range('E51').Value = '  +0.21%  '
